$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text even when it looks like a
# percentage (e.g. "76%"), which Excel would otherwise auto-convert to a
# numeric percentage. We stage the text (apostrophe-prefixed so it is
# stored as text) in a scratch cell far outside the used range, copy it,
# and paste-special (values only) into the target cell so the target
# keeps its own original style/number format untouched. The scratch cell
# is cleared afterwards so it leaves no trace in the sheet.
function Set-TextValue($cellRef, $val) {
    $ws.Range('ZZ1000').Value = "'" + $val
    $ws.Range('ZZ1000').Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
    $ws.Range('ZZ1000').Clear()
}

$ws.Range('E2').Value = '2026-02-12 21:18:28'
$ws.Range('E3').Value = '2026-02-12 21:18:30'
$ws.Range('E4').Value = '2026-02-12 21:18:33'
$ws.Range('J4').Value = '999.5 hPa'
$ws.Range('E5').Value = '2026-02-12 21:18:35'
$ws.Range('E6').Value = '2026-02-12 21:18:38'
$ws.Range('J6').Value = '999.3 hPa'
$ws.Range('E7').Value = '2026-02-12 21:18:40'
$ws.Range('J7').Value = '1002.0 hPa'
$ws.Range('O7').Value = '17.0 °C'
$ws.Range('E8').Value = '2026-02-12 21:18:43'
$ws.Range('J8').Value = '1001.4 hPa'
$ws.Range('K8').Value = '14.3 MJ/m2'
$ws.Range('E9').Value = '2026-02-12 21:18:46'
$ws.Range('O9').Value = '12.8 °C'
$ws.Range('E10').Value = '2026-02-12 21:18:48'
$ws.Range('K10').Value = '13.5 MJ/m2'
$ws.Range('O10').Value = '14.8 °C'
$ws.Range('E11').Value = '2026-02-12 21:18:51'
$ws.Range('O11').Value = '9.4 °C'
$ws.Range('E12').Value = '2026-02-12 21:18:53'
$ws.Range('E13').Value = '2026-02-12 21:18:55'
$ws.Range('J13').Value = '1002.0 hPa'
$ws.Range('E14').Value = '2026-02-12 21:18:58'
$ws.Range('E15').Value = '2026-02-12 21:19:01'
$ws.Range('E16').Value = '2026-02-12 21:19:03'
$ws.Range('E17').Value = '2026-02-12 21:19:06'
Set-TextValue 'H17' '76%'
$ws.Range('E18').Value = '2026-02-12 21:19:08'
$ws.Range('J18').Value = '999.7 hPa'
$ws.Range('E19').Value = '2026-02-12 21:19:11'
$ws.Range('E20').Value = '2026-02-12 21:19:13'
$ws.Range('O20').Value = '-3.5 °C'
$ws.Range('E21').Value = '2026-02-12 21:19:16'
$ws.Range('J21').Value = '1002.5 hPa'
$ws.Range('N21').Value = '5.1 °C 20:57 TU'
$ws.Range('O21').Value = '9.2 °C'
$ws.Range('E22').Value = '2026-02-12 21:19:19'
Set-TextValue 'H22' '77%'
$ws.Range('N22').Value = '-7.1 °C 20:44 TU'
$ws.Range('E23').Value = '2026-02-12 21:19:21'
$ws.Range('E24').Value = '2026-02-12 21:19:24'
$ws.Range('J24').Value = '1006.7 hPa'
$ws.Range('E25').Value = '2026-02-12 21:19:26'
$ws.Range('E26').Value = '2026-02-12 21:19:29'
$ws.Range('J26').Value = '999.0 hPa'
$ws.Range('O26').Value = '5.9 °C'
$ws.Range('E27').Value = '2026-02-12 21:19:32'
$ws.Range('E28').Value = '2026-02-12 21:19:34'
Set-TextValue 'H28' '39%'
$ws.Range('J28').Value = '999.2 hPa'
$ws.Range('N28').Value = '6.5 °C 20:56 TU'
$ws.Range('O28').Value = '13.9 °C'
$ws.Range('E29').Value = '2026-02-12 21:19:37'
Set-TextValue 'H29' '59%'
$ws.Range('K29').Value = '13.5 MJ/m2'
$ws.Range('N29').Value = '6.2 °C 20:49 TU'
$ws.Range('O29').Value = '14.4 °C'
$ws.Range('E30').Value = '2026-02-12 21:19:39'
$ws.Range('J30').Value = '999.6 hPa'
$ws.Range('N30').Value = '7.4 °C 20:46 TU'
$ws.Range('O30').Value = '12.1 °C'
$ws.Range('E31').Value = '2026-02-12 21:19:42'
$ws.Range('J31').Value = '999.0 hPa'
$ws.Range('E32').Value = '2026-02-12 21:19:45'
$ws.Range('E33').Value = '2026-02-12 21:19:47'
$ws.Range('J33').Value = '1001.6 hPa'
$ws.Range('E34').Value = '2026-02-12 21:19:50'
$ws.Range('E35').Value = '2026-02-12 21:19:52'
$ws.Range('E36').Value = '2026-02-12 21:19:55'
Set-TextValue 'H36' '59%'
$ws.Range('J36').Value = '999.9 hPa'
$ws.Range('E37').Value = '2026-02-12 21:19:58'
$ws.Range('J37').Value = '1000.5 hPa'
$ws.Range('N37').Value = '5.6 °C 20:54 TU'
$ws.Range('E38').Value = '2026-02-12 21:20:00'
$ws.Range('O38').Value = '15.8 °C'
$ws.Range('E39').Value = '2026-02-12 21:20:03'
$ws.Range('O39').Value = '-3.5 °C'
$ws.Range('E40').Value = '2026-02-12 21:20:05'
$ws.Range('J40').Value = '1003.3 hPa'
$ws.Range('O40').Value = '9.5 °C'
$ws.Range('E41').Value = '2026-02-12 21:20:08'
$ws.Range('J41').Value = '1005.7 hPa'
$ws.Range('E42').Value = '2026-02-12 21:20:10'
$ws.Range('N42').Value = '7.7 °C 20:59 TU'
$ws.Range('O42').Value = '14.1 °C'
$ws.Range('E43').Value = '2026-02-12 21:20:13'
$ws.Range('E44').Value = '2026-02-12 21:20:15'
$ws.Range('E45').Value = '2026-02-12 21:20:18'
$ws.Range('J45').Value = '1005.2 hPa'
$ws.Range('O45').Value = '6.9 °C'
$ws.Range('E46').Value = '2026-02-12 21:20:20'

$ws.Application.CutCopyMode = $false
